# Update "Código Efecto Hall" workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update header labels: change unit from (mV) to (V) ---
$ws.Range("F1").Value = "VL_-20mA (V)"
$ws.Range("G1").Value = "VL_-10mA (V)"
$ws.Range("H1").Value = "VL_10mA (V)"
$ws.Range("I1").Value = "VL_20mA (V)"
$ws.Range("J1").Value = "VL_err (V)"

# --- Update VL_err column values from 0.1 to 1E-3 (0.001) for rows 2-16 ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.001
}

# --- Update sheet view state: top-left cell and selected cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G11").Select()
